$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) cells whose new value looks like a plain decimal
# number (single '.') must be pre-formatted as Text, otherwise Excel's
# COM layer auto-converts the typed string into a numeric value (losing
# the original text-cell semantics used throughout this sheet). Cells
# whose new value contains more than one '.' (e.g. "66.217.95") are never
# auto-converted, so they don't need this treatment.
$textPriceRows = @(5,6,7,9,10,15,18,19,20,29,31,34,36,38,39,41,42,43,45,46)
foreach ($r in $textPriceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "66.217.95"
$ws.Range("E2").Value = "  -0.95%  "

$ws.Range("D3").Value = "3.275.16"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "586.21"
$ws.Range("E5").Value = "  +1.98%  "

$ws.Range("D6").Value = "178.96"
$ws.Range("E6").Value = "  -2.10%  "

$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  +4.17%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "0.126"
$ws.Range("E9").Value = "  -2.32%  "

$ws.Range("D10").Value = "6.74"
$ws.Range("E10").Value = "  +1.67%  "

$ws.Range("E11").Value = "  -0.61%  "

$ws.Range("D12").Value = "3.848.20"
$ws.Range("E12").Value = "  -1.34%  "

$ws.Range("E13").Value = "  -3.41%  "

$ws.Range("D14").Value = "66.207.24"
$ws.Range("E14").Value = "  -1.23%  "

$ws.Range("D15").Value = "26.33"
$ws.Range("E15").Value = "  -2.90%  "

$ws.Range("E16").Value = "  -1.62%  "

$ws.Range("D17").Value = "3.286.19"
$ws.Range("E17").Value = "  -1.71%  "

$ws.Range("D18").Value = "432.73"
$ws.Range("E18").Value = "  -1.57%  "

$ws.Range("D19").Value = "5.51"
$ws.Range("E19").Value = "  -2.43%  "

$ws.Range("D20").Value = "13.16"
$ws.Range("E20").Value = "  -2.74%  "

$ws.Range("E21").Value = "  -3.95%  "

$ws.Range("E22").Value = "  -2.95%  "

$ws.Range("E23").Value = "  +0.29%  "

$ws.Range("D24").Value = "3.422.77"
$ws.Range("E24").Value = "  -1.53%  "

$ws.Range("E25").Value = "  -1.05%  "

$ws.Range("E26").Value = "  +0.86%  "

$ws.Range("E27").Value = "  -4.88%  "

$ws.Range("E28").Value = "  -1.39%  "

$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("D31").Value = "22.25"
$ws.Range("E31").Value = "  -2.72%  "

$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("E33").Value = "  -2.43%  "

$ws.Range("D34").Value = "6.60"
$ws.Range("E34").Value = "  -2.47%  "

$ws.Range("E35").Value = "  -2.36%  "

$ws.Range("D36").Value = "157.80"
$ws.Range("E36").Value = "  -2.23%  "

$ws.Range("E37").Value = "  -4.74%  "

$ws.Range("D38").Value = "26.47"
$ws.Range("E38").Value = "  -3.70%  "

$ws.Range("D39").Value = "1.79"
$ws.Range("E39").Value = "  -2.76%  "

$ws.Range("D40").Value = "2.775.67"
$ws.Range("E40").Value = "  -1.49%  "

$ws.Range("D41").Value = "0.773"
$ws.Range("E41").Value = "  -1.81%  "

$ws.Range("D42").Value = "4.33"
$ws.Range("E42").Value = "  -2.49%  "

$ws.Range("D43").Value = "40.20"
$ws.Range("E43").Value = "  -0.27%  "

$ws.Range("E44").Value = "  -2.70%  "

$ws.Range("D45").Value = "0.0658"
$ws.Range("E45").Value = "  -1.77%  "

$ws.Range("D46").Value = "320.90"
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("E47").Value = "  -0.87%  "

$ws.Range("E48").Value = "  -4.43%  "

$ws.Range("E49").Value = "  -2.02%  "

$ws.Range("E50").Value = "  +3.14%  "

$ws.Range("E51").Value = "  +0.06%  "
